$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Add the new ammo row (row 44): ammo_og-7b / EX ---
# Bring over the number formats used by row 43 (D:K) so the new row matches
# the sheet's existing look (price/round, damage-per-rouble, pen-per-rouble,
# in-game power columns) before filling in the real values/formulas.
$ws.Range("D43:K43").Copy()
$ws.Range("D44:K44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A44").Value = "ammo_og-7b"
$ws.Range("B44").Value = "EX"
$ws.Range("C44").Value = 12490
$ws.Range("D44").Formula = "=C44/30"
$ws.Range("E44").Formula = "=K44/D44"
$ws.Range("F44").Formula = "=G44/D44*100"
$ws.Range("G44").Value = 0.37
$ws.Range("H44").Value = 3
$ws.Range("I44").Value = 3
$ws.Range("J44").Formula = "=I44*H44"
$ws.Range("K44").Formula = "=J44*Feuil2!`$B`$1"

# --- Restore the sheet's selection/scroll state recorded in the file ---
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D26").Select() | Out-Null
